$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parametrosInicio")
$ws.Range("B10").Value = "PARA BOOT MIGRACIONES SGV ENERO 2023 16.02.2023"
